$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 449.75
$ws.Range("I5").Value = 139.6
$ws.Range("J5").Value = 966.6667
$ws.Range("K5").Value = 139.6
$ws.Range("L5").Value = 966.6667
$ws.Range("M5").Value = -24.59999999999999
$ws.Range("N5").Value = -1196.6667
$ws.Range("H41").Value = 292
$ws.Range("J41").Value = 126
$ws.Range("L41").Value = 126
$ws.Range("N41").Value = -1006
$ws.Range("H64").Value = 5107.75
$ws.Range("I64").Value = 3722.111
$ws.Range("K64").Value = 3722.111
$ws.Range("M64").Value = -3474.111
$ws.Range("H67").Value = 5107.75
$ws.Range("I67").Value = 3722.111
$ws.Range("K67").Value = 3722.111
$ws.Range("M67").Value = -2864.111
$ws.Range("H125").Value = 4272.5713
$ws.Range("I125").Value = 4057.3333
$ws.Range("K125").Value = 36515.9997
$ws.Range("M125").Value = -34055.9997
$ws.Range("H138").Value = 1920.23
$ws.Range("I138").Value = 1459.3334
$ws.Range("J138").Value = 1965.8132
$ws.Range("K138").Value = 4378.0002
$ws.Range("L138").Value = 5897.4396
$ws.Range("M138").Value = 761.9997999999996
$ws.Range("N138").Value = -16177.4396

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3364.3704
$ws.Range("I45").Value = 926.55554
$ws.Range("K45").Value = 926.55554
$ws.Range("M45").Value = -549.55554
$ws.Range("H88").Value = 16078.143
$ws.Range("J88").Value = 16078.143
$ws.Range("L88").Value = 16078.143
$ws.Range("N88").Value = -16890.143
$ws.Range("H91").Value = 16078.143
$ws.Range("J91").Value = 16078.143
$ws.Range("L91").Value = 16078.143
$ws.Range("N91").Value = -18886.143
$ws.Range("H102").Value = 7867.385
$ws.Range("I102").Value = 7867.385
$ws.Range("K102").Value = 7867.385
$ws.Range("M102").Value = -6245.385

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2174.5557
$ws.Range("I20").Value = 2454.5715
$ws.Range("K20").Value = 2454.5715
$ws.Range("M20").Value = -2207.5715
$ws.Range("H22").Value = 2998
$ws.Range("J22").Value = 2998
$ws.Range("L22").Value = 2998
$ws.Range("N22").Value = -3344
$ws.Range("H86").Value = 2329.3333
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 2329.3333
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H92").Value = 34749.5
$ws.Range("J92").Value = 34749.5
$ws.Range("L92").Value = 34749.5
$ws.Range("N92").Value = -39741.5
$ws.Range("H94").Value = 750.9
$ws.Range("I94").Value = 439.66666
$ws.Range("K94").Value = 439.66666
$ws.Range("M94").Value = 11.33334000000002
$ws.Range("H107").Value = 1395
$ws.Range("I107").Value = 1201.7778
$ws.Range("J107").Value = 1742.8
$ws.Range("K107").Value = 1201.7778
$ws.Range("L107").Value = 1742.8
$ws.Range("M107").Value = 718.2221999999999
$ws.Range("N107").Value = -5582.8
$ws.Range("H130").Value = 90763
$ws.Range("J130").Value = 90763
$ws.Range("L130").Value = 90763
$ws.Range("N130").Value = -100803
$ws.Range("H132").Value = 81799.60000000001
$ws.Range("J132").Value = 81799.60000000001
$ws.Range("L132").Value = 81799.60000000001
$ws.Range("N132").Value = -91919.60000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 6284.2144
$ws.Range("I86").Value = 5579.5
$ws.Range("J86").Value = 6812.75
$ws.Range("K86").Value = 5579.5
$ws.Range("L86").Value = 6812.75
$ws.Range("M86").Value = -4456.5
$ws.Range("N86").Value = -9058.75
$ws.Range("H89").Value = 6284.2144
$ws.Range("I89").Value = 5579.5
$ws.Range("J89").Value = 6812.75
$ws.Range("K89").Value = 27897.5
$ws.Range("L89").Value = 34063.75
$ws.Range("M89").Value = -22281.5
$ws.Range("N89").Value = -45295.75
$ws.Range("H132").Value = 4311
$ws.Range("I132").Value = 4396.2856
$ws.Range("K132").Value = 13188.8568
$ws.Range("M132").Value = -10658.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 84
$ws.Range("I2").Value = 62.5
$ws.Range("J2").Value = 94.75
$ws.Range("K2").Value = 375
$ws.Range("L2").Value = 568.5
$ws.Range("M2").Value = -262
$ws.Range("N2").Value = -794.5
$ws.Range("H5").Value = 800.5925999999999
$ws.Range("I5").Value = 611.7895
$ws.Range("J5").Value = 1249
$ws.Range("K5").Value = 1835.3685
$ws.Range("L5").Value = 3747
$ws.Range("M5").Value = -1723.3685
$ws.Range("N5").Value = -3971
$ws.Range("H17").Value = 3421.5557
$ws.Range("J17").Value = 2949.3333
$ws.Range("L17").Value = 8847.999899999999
$ws.Range("N17").Value = -9185.999899999999
$ws.Range("H55").Value = 6058.3335
$ws.Range("J55").Value = 7655.143
$ws.Range("L55").Value = 22965.429
$ws.Range("N55").Value = -23319.429
$ws.Range("H68").Value = 1346.7
$ws.Range("I68").Value = 525
$ws.Range("J68").Value = 1552.125
$ws.Range("K68").Value = 1575
$ws.Range("L68").Value = 4656.375
$ws.Range("M68").Value = -764
$ws.Range("N68").Value = -6278.375
$ws.Range("H71").Value = 1346.7
$ws.Range("I71").Value = 525
$ws.Range("J71").Value = 1552.125
$ws.Range("K71").Value = 4725
$ws.Range("L71").Value = 13969.125
$ws.Range("M71").Value = -669
$ws.Range("N71").Value = -22081.125
$ws.Range("H122").Value = 2481.125
$ws.Range("J122").Value = 2942
$ws.Range("L122").Value = 26478
$ws.Range("N122").Value = -31378
$ws.Range("H132").Value = 6179.643
$ws.Range("I132").Value = 3144.6
$ws.Range("J132").Value = 7865.778
$ws.Range("K132").Value = 28301.4
$ws.Range("L132").Value = 70792.00200000001
$ws.Range("M132").Value = -25771.4
$ws.Range("N132").Value = -75852.00200000001
$ws.Range("H135").Value = 800.5925999999999
$ws.Range("I135").Value = 611.7895
$ws.Range("J135").Value = 1249
$ws.Range("K135").Value = 5506.1055
$ws.Range("L135").Value = 11241
$ws.Range("M135").Value = -2971.1055
$ws.Range("N135").Value = -16311

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 531.25
$ws.Range("I2").Value = 72
$ws.Range("K2").Value = 72
$ws.Range("M2").Value = 41
$ws.Range("H70").Value = 11875.5
$ws.Range("I70").Value = 15392.2
$ws.Range("J70").Value = 8358.799999999999
$ws.Range("K70").Value = 15392.2
$ws.Range("L70").Value = 8358.799999999999
$ws.Range("M70").Value = -15122.2
$ws.Range("N70").Value = -8898.799999999999
$ws.Range("H73").Value = 11875.5
$ws.Range("I73").Value = 15392.2
$ws.Range("J73").Value = 8358.799999999999
$ws.Range("K73").Value = 15392.2
$ws.Range("L73").Value = 8358.799999999999
$ws.Range("M73").Value = -14456.2
$ws.Range("N73").Value = -10230.8
$ws.Range("H126").Value = 12495.75
$ws.Range("I126").Value = 9999.75
$ws.Range("K126").Value = 29999.25
$ws.Range("M126").Value = -27529.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3099.6667
$ws.Range("J40").Value = 4299
$ws.Range("L40").Value = 4299
$ws.Range("N40").Value = -4571
$ws.Range("H55").Value = 753.7
$ws.Range("I55").Value = 172.375
$ws.Range("J55").Value = 1418.0714
$ws.Range("K55").Value = 172.375
$ws.Range("L55").Value = 1418.0714
$ws.Range("M55").Value = 0.625
$ws.Range("N55").Value = -1764.0714

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H126").Value = 2363.4736
$ws.Range("I126").Value = 1887.6666
$ws.Range("K126").Value = 5662.9998
$ws.Range("M126").Value = -3192.9998
$ws.Range("H132").Value = 3003.4443
$ws.Range("J132").Value = 7000
$ws.Range("L132").Value = 21000
$ws.Range("N132").Value = -26060
